$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44202
$ws.Range("I2").Value = "Segunda"
$ws.Range("J2").Value = 1300
$ws.Range("K2").Value = 230
$ws.Range("L2").Value = 250
$ws.Range("M2").Value = 240
$ws.Range("P2").Value = 240
$ws.Range("D3").Value = 44231
$ws.Range("J3").Value = 200
$ws.Range("K3").Value = 180
$ws.Range("L3").Value = 200
$ws.Range("M3").Value = 190
$ws.Range("O3").Value = "Región de Arica y Parinacota"
$ws.Range("P3").Value = 190
$ws.Range("D4").Value = 44224
$ws.Range("I4").Value = "Segunda"
$ws.Range("K4").Value = 230
$ws.Range("L4").Value = 250
$ws.Range("M4").Value = 240
$ws.Range("P4").Value = 240
$ws.Range("D5").Value = 44224
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 200
$ws.Range("L5").Value = 230
$ws.Range("M5").Value = 215
$ws.Range("O5").Value = "Región de Arica y Parinacota"
$ws.Range("P5").Value = 215
$ws.Range("D6").Value = 44229
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 1200
$ws.Range("K6").Value = 230
$ws.Range("M6").Value = 240
$ws.Range("P6").Value = 240
$ws.Range("D7").Value = 44253
$ws.Range("I7").Value = "Segunda"
$ws.Range("J7").Value = 1200
$ws.Range("K7").Value = 270
$ws.Range("L7").Value = 280
$ws.Range("M7").Value = 275
$ws.Range("P7").Value = 275
$ws.Range("D8").Value = 44214
$ws.Range("J8").Value = 1200
$ws.Range("K8").Value = 400
$ws.Range("L8").Value = 450
$ws.Range("M8").Value = 425
$ws.Range("P8").Value = 425
$ws.Range("D9").Value = 44217
$ws.Range("J9").Value = 1600
$ws.Range("K9").Value = 300
$ws.Range("L9").Value = 350
$ws.Range("M9").Value = 325
$ws.Range("P9").Value = 325
$ws.Range("D10").Value = 44172
$ws.Range("J10").Value = 1600
$ws.Range("K10").Value = 400
$ws.Range("L10").Value = 420
$ws.Range("M10").Value = 410
$ws.Range("P10").Value = 410
$ws.Range("D11").Value = 44201
$ws.Range("I11").Value = "Segunda"
$ws.Range("J11").Value = 1800
$ws.Range("K11").Value = 250
$ws.Range("L11").Value = 270
$ws.Range("M11").Value = 260
$ws.Range("P11").Value = 260
$ws.Range("D12").Value = 44301
$ws.Range("J12").Value = 900
$ws.Range("K12").Value = 280
$ws.Range("L12").Value = 300
$ws.Range("M12").Value = 290
$ws.Range("O12").Value = "Perú"
$ws.Range("P12").Value = 290
$ws.Range("D13").Value = 44251
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 1200
$ws.Range("K13").Value = 250
$ws.Range("L13").Value = 280
$ws.Range("M13").Value = 265
$ws.Range("P13").Value = 265
$ws.Range("D14").Value = 44243
$ws.Range("I14").Value = "Primera"
$ws.Range("K14").Value = 300
$ws.Range("L14").Value = 320
$ws.Range("M14").Value = 310
$ws.Range("P14").Value = 310
$ws.Range("D15").Value = 44243
$ws.Range("J15").Value = 800
$ws.Range("K15").Value = 300
$ws.Range("L15").Value = 320
$ws.Range("M15").Value = 310
$ws.Range("O15").Value = "Perú"
$ws.Range("P15").Value = 310
$ws.Range("D16").Value = 44166
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 1700
$ws.Range("K16").Value = 500
$ws.Range("L16").Value = 530
$ws.Range("M16").Value = 515
$ws.Range("P16").Value = 515
$ws.Range("D18").Value = 44160
$ws.Range("J18").Value = 2000
$ws.Range("K18").Value = 500
$ws.Range("L18").Value = 550
$ws.Range("M18").Value = 525
$ws.Range("P18").Value = 525
$ws.Range("D19").Value = 44175
$ws.Range("J19").Value = 1200
$ws.Range("K19").Value = 400
$ws.Range("L19").Value = 430
$ws.Range("M19").Value = 415
$ws.Range("P19").Value = 415
$ws.Range("D20").Value = 44162
$ws.Range("J20").Value = 900
$ws.Range("K20").Value = 500
$ws.Range("L20").Value = 550
$ws.Range("M20").Value = 525
$ws.Range("P20").Value = 525
$ws.Range("D21").Value = 44162
$ws.Range("I21").Value = "Segunda"
$ws.Range("K21").Value = 500
$ws.Range("L21").Value = 550
$ws.Range("M21").Value = 525
$ws.Range("P21").Value = 525
$ws.Range("D22").Value = 44176
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 1300
$ws.Range("K22").Value = 350
$ws.Range("M22").Value = 375
$ws.Range("P22").Value = 375

"Done applying changes"